$wb = $excel.ActiveWorkbook

# --- Create "Norway" sheet as a copy of "Hungary" (placed right after it) ---
$hungary = $wb.Worksheets.Item("Hungary")
$hungary.Copy($null, $hungary)
$norway = $wb.Worksheets.Item($hungary.Index + 1)
$norway.Name = "Norway"

# Update the market-specific values
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-2931/T3062"

# Column D width becomes a fixed custom width of 25 (was bestFit 8.44)
$norway.Columns.Item(4).ColumnWidth = 24.166666666666668

# Rows 3-5 lose their explicit 28.8 height override (back to sheet default)
$norway.Rows.Item(3).EntireRow.AutoFit()
$norway.Rows.Item(4).EntireRow.AutoFit()
$norway.Rows.Item(5).EntireRow.AutoFit()

# Whole sheet is selected on this tab
$norway.Cells.Select() | Out-Null

# --- Create "Poland" sheet as a copy of "Norway" (placed right after it) ---
$norway.Copy($null, $norway)
$poland = $wb.Worksheets.Item($norway.Index + 1)
$poland.Name = "Poland"

$poland.Range("B2").Value = "Poland Market"
$poland.Range("B4").Value = "NGC-2920/T3105"

$poland.Columns.Item(4).ColumnWidth = 24.166666666666668
$poland.Rows.Item(3).EntireRow.AutoFit()
$poland.Rows.Item(4).EntireRow.AutoFit()
$poland.Rows.Item(5).EntireRow.AutoFit()
$poland.Cells.Select() | Out-Null

# The newly added "Norway" tab is the active/visible one
$norway.Select()
